$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tendencia")

# --- Column E: new "label" column with tendency formula ---
$ws.Cells.Item(1, 5).Value = "label"

for ($r = 2; $r -le 47; $r++) {
    $formula = '=IF(D' + $r + '=1, "Derecha", IF(D' + $r + '=-1, "Izquierda", "Sin tendencia"))'
    $ws.Cells.Item($r, 5).Formula = $formula
}

# --- Column A: refresh election labels (strip trailing space / reorder) ---
for ($r = 9; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = "2017_presidencial_1v"
}
for ($r = 19; $r -le 22; $r++) {
    $ws.Cells.Item($r, 1).Value = "2017_presidencial_2v"
}
for ($r = 35; $r -le 43; $r++) {
    $ws.Cells.Item($r, 1).Value = "2021_presidencial_1v"
}
for ($r = 44; $r -le 47; $r++) {
    $ws.Cells.Item($r, 1).Value = "2021_presidencial_2v"
}

# --- Column E width ---
$ws.Columns.Item(5).ColumnWidth = 12.3

# --- View: active cell / selection on the tendencia sheet ---
$ws.Activate()
$ws.Range("C28").Select()

# --- Workbook view: first visible tab scrolled to sheet 3 (tendencia) ---
$win = $excel.ActiveWindow
$win.ScrollWorkbookTabs(3)
